$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at U (shifts existing U:W -> V:X, i.e.
# insuf_venosa_periferica/obesidade/SAOS move one column right) and
# populate it with the new "doenca_renal_cronica" variable.
$ws.Columns("U:U").Insert()

# Header
$ws.Range("U1").Value = "doenca_renal_cronica"

# Data rows: default everybody to 0 ...
$ws.Range("U2:U64").Value = 0

# ... except the patient with id = 22 (spreadsheet row 23), who has the
# condition.
$ws.Range("U23").Value = 1
